$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: insert Jan_2026, shift Dec_2025/Nov_2025, drop Oct_2025
$ws.Cells.Item(1,4).Value = "Jan_2026"
$ws.Cells.Item(1,5).Value = "Dec_2025"
$ws.Cells.Item(1,6).Value = "Nov_2025"

# Row 2
$ws.Cells.Item(2,1).Value = "INE795G01014"
$ws.Cells.Item(2,2).Value = "HDFC Life Insurance Co Ltd"
$ws.Cells.Item(2,4).Value = 9.956272999999999
$ws.Cells.Item(2,5).Value = 6.706992
$ws.Cells.Item(2,6).Value = 6.598957
$ws.Cells.Item(2,7).Value = 3.249281
$ws.Cells.Item(2,8).Value = 3.357315999999999

# Row 3
$ws.Cells.Item(3,1).Value = "INE775A01035"
$ws.Cells.Item(3,2).Value = "Samvardhana Motherson International Ltd"
$ws.Cells.Item(3,4).Value = 9.790407999999999
$ws.Cells.Item(3,5).Value = 9.361610000000001
$ws.Cells.Item(3,6).Value = 6.876986
$ws.Cells.Item(3,7).Value = 0.4287979999999987
$ws.Cells.Item(3,8).Value = 2.913422

# Row 4
$ws.Cells.Item(4,1).Value = "INE040A01034"
$ws.Cells.Item(4,2).Value = "HDFC Bank Limited"
$ws.Cells.Item(4,4).Value = 9.502613
$ws.Cells.Item(4,5).Value = 2.93099
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 6.571623000000001
$ws.Cells.Item(4,8).Value = 9.502613

# Row 5
$ws.Cells.Item(5,1).Value = "INE237A01036"
$ws.Cells.Item(5,2).Value = "Kotak Mahindra Bank Limited"
$ws.Cells.Item(5,4).Value = 8.638477999999999
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 0
$ws.Cells.Item(5,7).Value = 8.638477999999999
$ws.Cells.Item(5,8).Value = 8.638477999999999

# Row 6
$ws.Cells.Item(6,1).Value = "INE522F01014"
$ws.Cells.Item(6,2).Value = "Coal India Ltd"
$ws.Cells.Item(6,4).Value = 7.5848
$ws.Cells.Item(6,5).Value = 0
$ws.Cells.Item(6,6).Value = 0
$ws.Cells.Item(6,7).Value = 7.5848
$ws.Cells.Item(6,8).Value = 7.5848

# Row 7
$ws.Cells.Item(7,1).Value = "INE090A01021"
$ws.Cells.Item(7,2).Value = "ICICI Bank Limited"
$ws.Cells.Item(7,4).Value = 6.567132
$ws.Cells.Item(7,5).Value = 3.016362
$ws.Cells.Item(7,6).Value = 3.011186
$ws.Cells.Item(7,7).Value = 3.55077
$ws.Cells.Item(7,8).Value = 3.555946

# Row 8
$ws.Cells.Item(8,1).Value = "INE406A01037"
$ws.Cells.Item(8,2).Value = "Aurobindo Pharma Limited"
$ws.Cells.Item(8,4).Value = 6.310856
$ws.Cells.Item(8,5).Value = 5.564083
$ws.Cells.Item(8,6).Value = 5.569359
$ws.Cells.Item(8,7).Value = 0.7467730000000001
$ws.Cells.Item(8,8).Value = 0.7414969999999999

# Row 9
$ws.Cells.Item(9,1).Value = "INE781S01027"
$ws.Cells.Item(9,2).Value = "Ventive Hospitality Limited"
$ws.Cells.Item(9,4).Value = 6.171473
$ws.Cells.Item(9,5).Value = 5.68026
$ws.Cells.Item(9,6).Value = 5.233346
$ws.Cells.Item(9,7).Value = 0.4912130000000001
$ws.Cells.Item(9,8).Value = 0.9381269999999997

# Row 10
$ws.Cells.Item(10,1).Value = "INE917I01010"
$ws.Cells.Item(10,2).Value = "Bajaj Auto Limited"
$ws.Cells.Item(10,4).Value = 4.906354
$ws.Cells.Item(10,5).Value = 4.298993
$ws.Cells.Item(10,6).Value = 4.030076
$ws.Cells.Item(10,7).Value = 0.607361
$ws.Cells.Item(10,8).Value = 0.8762780000000001

# Row 11
$ws.Cells.Item(11,1).Value = "INE364U01010"
$ws.Cells.Item(11,2).Value = "Adani Green Energy Limited"
$ws.Cells.Item(11,4).Value = 4.122146
$ws.Cells.Item(11,5).Value = 4.418953
$ws.Cells.Item(11,6).Value = 4.404663
$ws.Cells.Item(11,7).Value = -0.2968070000000003
$ws.Cells.Item(11,8).Value = -0.2825170000000004

# Row 12
$ws.Cells.Item(12,1).Value = "INE0BYP01024"
$ws.Cells.Item(12,2).Value = "OSWAL PUMPS LIMITED"
$ws.Cells.Item(12,4).Value = 1.849712
$ws.Cells.Item(12,5).Value = 2.216061
$ws.Cells.Item(12,6).Value = 2.237151
$ws.Cells.Item(12,7).Value = -0.3663489999999998
$ws.Cells.Item(12,8).Value = -0.3874389999999999

# Row 13
$ws.Cells.Item(13,1).Value = "INE0Z8Z23013"
$ws.Cells.Item(13,2).Value = "Capital Infra Trust InvIT"
$ws.Cells.Item(13,4).Value = 0.08158700000000001
$ws.Cells.Item(13,5).Value = 0.07441300000000001
$ws.Cells.Item(13,6).Value = 0
$ws.Cells.Item(13,7).Value = 0.007174
$ws.Cells.Item(13,8).Value = 0.08158700000000001

# Row 14
$ws.Cells.Item(14,1).Value = "INE271C01023"
$ws.Cells.Item(14,2).Value = "DLF Limited"
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 2.922047
$ws.Cells.Item(14,6).Value = 2.969165
$ws.Cells.Item(14,7).Value = -2.922047
$ws.Cells.Item(14,8).Value = -2.969165

# Row 15
$ws.Cells.Item(15,1).Value = "INE245A01021"
$ws.Cells.Item(15,2).Value = "Tata Power Company Limited"
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 0
$ws.Cells.Item(15,6).Value = 6.408252
$ws.Cells.Item(15,7).Value = 0
$ws.Cells.Item(15,8).Value = -6.408252

# Row 16
$ws.Cells.Item(16,1).Value = "INE484J01027"
$ws.Cells.Item(16,2).Value = "Godrej Properties Limited"
$ws.Cells.Item(16,4).Value = 0
$ws.Cells.Item(16,5).Value = 2.567573
$ws.Cells.Item(16,6).Value = 2.614716
$ws.Cells.Item(16,7).Value = -2.567573
$ws.Cells.Item(16,8).Value = -2.614716

# Row 17
$ws.Cells.Item(17,1).Value = "INE237A01028"
$ws.Cells.Item(17,2).Value = "Kotak Mahindra Bank Limited"
$ws.Cells.Item(17,4).Value = 0
$ws.Cells.Item(17,5).Value = 8.389315
$ws.Cells.Item(17,6).Value = 4.946473
$ws.Cells.Item(17,7).Value = -8.389315
$ws.Cells.Item(17,8).Value = -4.946473

# Row 18
$ws.Cells.Item(18,1).Value = "INE154A01025"
$ws.Cells.Item(18,2).Value = "ITC Limited"
$ws.Cells.Item(18,4).Value = 0
$ws.Cells.Item(18,5).Value = 8.106636999999999
$ws.Cells.Item(18,6).Value = 7.849532
$ws.Cells.Item(18,7).Value = -8.106636999999999
$ws.Cells.Item(18,8).Value = -7.849532

# Row 19
$ws.Cells.Item(19,1).Value = "INE296A01032"
$ws.Cells.Item(19,2).Value = "Bajaj Finance Limited"
$ws.Cells.Item(19,4).Value = 0
$ws.Cells.Item(19,5).Value = 1.781996
$ws.Cells.Item(19,6).Value = 3.300623
$ws.Cells.Item(19,7).Value = -1.781996
$ws.Cells.Item(19,8).Value = -3.300623
